# Generate Report for Handoff
# Updates status from "In Translation" to "Ready for handoff" and refreshes
# the related timestamps on the Overview, zh-cn and de-de sheets, then
# widens the affected Status columns to fit the new text.

$wb = $excel.ActiveWorkbook

# Target column width (OOXML "width" attribute) is 17.2159881591797 chars;
# the COM ColumnWidth setter only accepts whole-pixel character widths, so
# 16.333333333333332 (== 98/6) is the closest representable value and is
# what the interop layer rounds to.
$statusColWidth = 16.333333333333332

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-09-01 05:03:18"
$wsOverview.Columns.Item(5).ColumnWidth = $statusColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $statusColWidth

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-09-01 05:03:14"
$wsZhCn.Columns.Item(3).ColumnWidth = $statusColWidth

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-09-01 05:03:18"
$wsDeDe.Columns.Item(3).ColumnWidth = $statusColWidth
